# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# The "K" column (column G) values were recalculated upstream and need to be
# rewritten into the worksheet for rows 2-36 (data rows, row 1 is the header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 6
    4  = 10
    5  = 5
    6  = 7
    7  = 3
    8  = 7
    9  = 7
    10 = 4
    11 = 7
    12 = 10
    13 = 7
    14 = 10
    15 = 9
    16 = 5
    17 = 9
    18 = 7
    19 = 4
    20 = 3
    21 = 6
    22 = 6
    23 = 9
    24 = 11
    25 = 4
    26 = 7
    27 = 7
    28 = 3
    29 = 12
    30 = 9
    31 = 6
    32 = 4
    33 = 5
    34 = 5
    35 = 4
    36 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
